# Secret Santa workbook: replace the real participant email addresses with
# a single dummy placeholder address, drop the per-cell mailto: hyperlinks
# that pointed at the old addresses, widen the email column so the longer
# placeholder text is readable, and move the active selection back to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing mailto hyperlinks (B2, B5, B6) entirely - we no longer
# want live links on the dummy address.
$ws.Hyperlinks.Delete()

# Replace every participant's email with the dummy placeholder address.
$ws.Range("B2").Value = "example@mail.com"
$ws.Range("B3").Value = "example@mail.com"
$ws.Range("B4").Value = "example@mail.com"
$ws.Range("B5").Value = "example@mail.com"
$ws.Range("B6").Value = "example@mail.com"

# Keep the email column styled like the old hyperlink text (underlined /
# hyperlink theme color) even though it's no longer a clickable link.
$ws.Range("B2:B6").Style = "Hyperlink"

# Widen column B so "example@mail.com" has breathing room.
$ws.Range("B1").ColumnWidth = 36.35

# Move the selection to B2 (first email cell) instead of the old B10.
$ws.Range("B2").Select() | Out-Null
